$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.813.16'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '1.628.65'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").Value = "'215.44"
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("D6").Value = "'0.5059"
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("D7").Value = "'1.005"
$ws.Range("D8").Value = "'0.06453"
$ws.Range("E8").Value = '  +1.39%  '
$ws.Range("D9").Value = "'0.2577"
$ws.Range("E9").Value = '  -0.25%  '
$ws.Range("D10").Value = "'19.38"
$ws.Range("E10").Value = '  -2.53%  '
$ws.Range("D11").Value = "'0.07798"
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").Value = "'4.255"
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").Value = '1.629.70'
$ws.Range("E13").Value = '  -0.48%  '
$ws.Range("D14").Value = '1.853.28'
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").Value = "'0.5568"
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("D16").Value = "'63.14"
$ws.Range("E16").Value = '  -1.89%  '
$ws.Range("D17").Value = '0.0₅7536'
$ws.Range("E17").Value = '  -2.85%  '
$ws.Range("D18").Value = '25.824.07'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("E19").Value = '  +0.28%  '
$ws.Range("D20").Value = "'193.20"
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("D21").Value = "'4.293"
$ws.Range("E21").Value = '  -3.32%  '
$ws.Range("D22").Value = "'9.796"
$ws.Range("E22").Value = '  -1.50%  '
$ws.Range("D23").Value = "'6.000"
$ws.Range("E23").Value = '  -1.68%  '
$ws.Range("E25").Value = '  -4.19%  '
$ws.Range("D26").Value = "'140.23"
$ws.Range("E26").Value = '  -2.53%  '
$ws.Range("D27").Value = "'0.1260"
$ws.Range("E27").Value = '  +1.65%  '
$ws.Range("D28").Value = "'6.714"
$ws.Range("E28").Value = '  -2.37%  '
$ws.Range("D29").Value = "'15.38"
$ws.Range("E29").Value = '  -1.76%  '
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").Value = "'0.04854"
$ws.Range("E31").Value = '  -0.67%  '
$ws.Range("D32").Value = "'3.271"
$ws.Range("E32").Value = '  -0.26%  '
$ws.Range("D33").Value = "'3.175"
$ws.Range("E33").Value = '  -1.61%  '
$ws.Range("D34").Value = "'1.552"
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D35").Value = "'2.380"
$ws.Range("E35").Value = '  +0.22%  '
$ws.Range("D36").Value = "'0.8925"
$ws.Range("E36").Value = '  -2.61%  '
$ws.Range("D37").Value = "'2.568"
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '1.133.34'
$ws.Range("E38").Value = '  +3.99%  '
$ws.Range("D39").Value = "'0.5455"
$ws.Range("E39").Value = '  -1.87%  '
$ws.Range("E40").Value = '  -1.16%  '
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = "'5.559"
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("D43").Value = "'0.7952"
$ws.Range("E43").Value = '  -1.30%  '
$ws.Range("D44").Value = "'97.20"
$ws.Range("E44").Value = '  -2.04%  '
$ws.Range("D45").Value = '1.781.76'
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("E46").Value = '  -7.64%  '
$ws.Range("D47").Value = "'0.4439"
$ws.Range("E47").Value = '  -2.15%  '
$ws.Range("D48").Value = "'54.94"
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").Value = "'0.05057"
$ws.Range("E49").Value = '  -3.14%  '
$ws.Range("D50").Value = "'7.604"
$ws.Range("E50").Value = '  +0.54%  '
$ws.Range("E51").Value = '  -0.15%  '
